# Apply trade #19 to the workbook: update summary stats and append the new
# trade row to the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = 0.34
$summary.Range("B6").Value = 19
$summary.Range("B9").Value = 52.63

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D6").Value = 19
$status.Range("G6").Value = 52.63

# --- New trade row data (row 20 on both "All Trades" and "MarketMaking") ---
$newRow = @(19, "2026-02-17", "23:54:31", "MarketMaking", "DOWN", 0.01, 0.01, "CLOSED", 0, 0, 100.32, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.11)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($col = 1; $col -le $newRow.Length; $col++) {
        # Column B holds a date-like string ("2026-02-17"); force it to be
        # stored as text so Excel doesn't auto-convert it to a date serial.
        if ($col -eq 2) {
            $ws.Cells.Item(20, $col).NumberFormat = "@"
        }
        $ws.Cells.Item(20, $col).Value = $newRow[$col - 1]
    }
}
